$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 7.781484444871653
$ws.Cells.Item(2, 4).Value = 8.994834548131863
$ws.Cells.Item(2, 5).Value = 13.03036528579197
$ws.Cells.Item(2, 6).Value = 32.43215308875988
$ws.Cells.Item(2, 7).Value = 3.629463912638744
$ws.Cells.Item(2, 10).Value = 9.683785215320425
$ws.Cells.Item(2, 13).Value = 20.46657744912371
$ws.Cells.Item(2, 14).Value = 17.86441640311219
$ws.Cells.Item(2, 15).Value = 24.53295025528709

# Row 3
$ws.Cells.Item(3, 2).Value = 7.707759290682338
$ws.Cells.Item(3, 4).Value = 9.020965173569156
$ws.Cells.Item(3, 5).Value = 13.09312267452674
$ws.Cells.Item(3, 6).Value = 32.28043253589874
$ws.Cells.Item(3, 7).Value = 3.632480853757775
$ws.Cells.Item(3, 10).Value = 9.723651660560115
$ws.Cells.Item(3, 13).Value = 19.93841419806747
$ws.Cells.Item(3, 14).Value = 17.79118311800559
$ws.Cells.Item(3, 15).Value = 24.41363490575228

# Row 4
$ws.Cells.Item(4, 2).Value = 7.663923939527042
$ws.Cells.Item(4, 4).Value = 9.037980542963863
$ws.Cells.Item(4, 5).Value = 13.13371433053429
$ws.Cells.Item(4, 6).Value = 32.19685971458989
$ws.Cells.Item(4, 7).Value = 3.634431784606534
$ws.Cells.Item(4, 10).Value = 9.749369017003557
$ws.Cells.Item(4, 13).Value = 19.60812315943738
$ws.Cells.Item(4, 14).Value = 17.74853499733086
$ws.Cells.Item(4, 15).Value = 24.34729758947682

# Row 5
$ws.Cells.Item(5, 2).Value = 7.646441234542565
$ws.Cells.Item(5, 4).Value = 9.04515904333485
$ws.Cells.Item(5, 5).Value = 13.15077431958728
$ws.Cells.Item(5, 6).Value = 32.16523803357006
$ws.Cells.Item(5, 7).Value = 3.635251664443437
$ws.Cells.Item(5, 10).Value = 9.760161453040576
$ws.Cells.Item(5, 13).Value = 19.47219856090954
$ws.Cells.Item(5, 14).Value = 17.73175216593198
$ws.Cells.Item(5, 15).Value = 24.32202543387413

# Row 6
$ws.Cells.Item(6, 2).Value = 7.643561812598325
$ws.Cells.Item(6, 4).Value = 9.046365813429153
$ws.Cells.Item(6, 5).Value = 13.15363846617399
$ws.Cells.Item(6, 6).Value = 32.16013499187181
$ws.Cells.Item(6, 7).Value = 3.635389308881889
$ws.Cells.Item(6, 10).Value = 9.761972420806355
$ws.Cells.Item(6, 13).Value = 19.44955360869221
$ws.Cells.Item(6, 14).Value = 17.72900181889754
$ws.Cells.Item(6, 15).Value = 24.31793589263802

# Row 7
$ws.Cells.Item(7, 2).Value = 7.663686594319385
$ws.Cells.Item(7, 4).Value = 9.038076363792198
$ws.Cells.Item(7, 5).Value = 13.13394230659842
$ws.Cells.Item(7, 6).Value = 32.19642336475315
$ws.Cells.Item(7, 7).Value = 3.634442741022303
$ws.Cells.Item(7, 10).Value = 9.749513301603823
$ws.Cells.Item(7, 13).Value = 19.6062951662267
$ws.Cells.Item(7, 14).Value = 17.7483062245642
$ws.Cells.Item(7, 15).Value = 24.34694960647707

# Row 8
$ws.Cells.Item(8, 2).Value = 7.755779213740088
$ws.Cells.Item(8, 4).Value = 9.003643096084458
$ws.Cells.Item(8, 5).Value = 13.05157740197411
$ws.Cells.Item(8, 6).Value = 32.37786674584754
$ws.Cells.Item(8, 7).Value = 3.630483759113902
$ws.Cells.Item(8, 10).Value = 9.697274451983388
$ws.Cells.Item(8, 13).Value = 20.28581585241087
$ws.Cells.Item(8, 14).Value = 17.83869308761492
$ws.Cells.Item(8, 15).Value = 24.49038637024339

# Row 9
$ws.Cells.Item(9, 2).Value = 7.946809038784773
$ws.Cells.Item(9, 4).Value = 8.943805448393343
$ws.Cells.Item(9, 5).Value = 12.90634774754657
$ws.Cells.Item(9, 6).Value = 32.80844874168471
$ws.Cells.Item(9, 7).Value = 3.623497939759467
$ws.Cells.Item(9, 10).Value = 9.604629873174368
$ws.Cells.Item(9, 13).Value = 21.56343011216754
$ws.Cells.Item(9, 14).Value = 18.03369547089514
$ws.Cells.Item(9, 15).Value = 24.82558393141726

# Row 10
$ws.Cells.Item(10, 2).Value = 8.092246164162297
$ws.Cells.Item(10, 4).Value = 8.904501688817589
$ws.Cells.Item(10, 5).Value = 12.80951405906424
$ws.Cells.Item(10, 6).Value = 33.16842505960729
$ws.Cells.Item(10, 7).Value = 3.618834059199366
$ws.Cells.Item(10, 10).Value = 9.542483458080284
$ws.Cells.Item(10, 13).Value = 22.45957790523367
$ws.Cells.Item(10, 14).Value = 18.18690673454109
$ws.Cells.Item(10, 15).Value = 25.10319202247406

# Row 11
$ws.Cells.Item(11, 2).Value = 8.159237429276244
$ws.Cells.Item(11, 4).Value = 8.88762748556265
$ws.Cells.Item(11, 5).Value = 12.76759122020191
$ws.Cells.Item(11, 6).Value = 33.34117341947903
$ws.Cells.Item(11, 7).Value = 3.616812914347657
$ws.Cells.Item(11, 10).Value = 9.515485763707495
$ws.Cells.Item(11, 13).Value = 22.85634917210396
$ws.Cells.Item(11, 14).Value = 18.25856875489477
$ws.Cells.Item(11, 15).Value = 25.23591115472478

# Row 12
$ws.Cells.Item(12, 2).Value = 8.184701560118098
$ws.Cells.Item(12, 4).Value = 8.881381820431553
$ws.Cells.Item(12, 5).Value = 12.75202111631603
$ws.Cells.Item(12, 6).Value = 33.40783629554103
$ws.Cells.Item(12, 7).Value = 3.616061917952958
$ws.Cells.Item(12, 10).Value = 9.505444687800093
$ws.Cells.Item(12, 13).Value = 23.00490572930469
$ws.Cells.Item(12, 14).Value = 18.28597070760701
$ws.Cells.Item(12, 15).Value = 25.28705862299429

# Row 13
$ws.Cells.Item(13, 2).Value = 8.179213507700805
$ws.Cells.Item(13, 4).Value = 8.882720526664707
$ws.Cells.Item(13, 5).Value = 12.75536085593248
$ws.Cells.Item(13, 6).Value = 33.39342456689306
$ws.Cells.Item(13, 7).Value = 3.616223020762812
$ws.Cells.Item(13, 10).Value = 9.507599113364806
$ws.Cells.Item(13, 13).Value = 22.97298856020934
$ws.Cells.Item(13, 14).Value = 18.28005769871777
$ws.Cells.Item(13, 15).Value = 25.27600412074786

# Row 14
$ws.Cells.Item(14, 2).Value = 8.161330589072451
$ws.Cells.Item(14, 4).Value = 8.887110762210533
$ws.Cells.Item(14, 5).Value = 12.76630414654986
$ws.Cells.Item(14, 6).Value = 33.34663310156601
$ws.Cells.Item(14, 7).Value = 3.616750841943373
$ws.Cells.Item(14, 10).Value = 9.514656027423282
$ws.Cells.Item(14, 13).Value = 22.86860560352266
$ws.Cells.Item(14, 14).Value = 18.26081791851106
$ws.Cells.Item(14, 15).Value = 25.24010146887921

# Row 15
$ws.Cells.Item(15, 2).Value = 8.150388590166937
$ws.Cells.Item(15, 4).Value = 8.889818681580229
$ws.Cells.Item(15, 5).Value = 12.77304694812049
$ws.Cells.Item(15, 6).Value = 33.31813292779817
$ws.Cells.Item(15, 7).Value = 3.617076016299302
$ws.Cells.Item(15, 10).Value = 9.519002319928024
$ws.Cells.Item(15, 13).Value = 22.80444411899551
$ws.Cells.Item(15, 14).Value = 18.2490669960556
$ws.Cells.Item(15, 15).Value = 25.2182248243805

# Row 16
$ws.Cells.Item(16, 2).Value = 8.087882838751103
$ws.Cells.Item(16, 4).Value = 8.905624658178235
$ws.Cells.Item(16, 5).Value = 12.81229655820636
$ws.Cells.Item(16, 6).Value = 33.15731261190308
$ws.Cells.Item(16, 7).Value = 3.618968160679084
$ws.Cells.Item(16, 10).Value = 9.544273377806372
$ws.Cells.Item(16, 13).Value = 22.43341778599028
$ws.Cells.Item(16, 14).Value = 18.18226151173332
$ws.Cells.Item(16, 15).Value = 25.09464498359441

# Row 17
$ws.Cells.Item(17, 2).Value = 8.049733695389826
$ws.Cells.Item(17, 4).Value = 8.915578352038258
$ws.Cells.Item(17, 5).Value = 12.83691924135268
$ws.Cells.Item(17, 6).Value = 33.06092582128262
$ws.Cells.Item(17, 7).Value = 3.620154606745882
$ws.Cells.Item(17, 10).Value = 9.560101907783672
$ws.Cells.Item(17, 13).Value = 22.20292245570065
$ws.Cells.Item(17, 14).Value = 18.14176954127647
$ws.Cells.Item(17, 15).Value = 25.02045592167496

# Row 18
$ws.Cells.Item(18, 2).Value = 8.027870957161024
$ws.Cells.Item(18, 4).Value = 8.92139809378069
$ws.Cells.Item(18, 5).Value = 12.85128183081259
$ws.Cells.Item(18, 6).Value = 33.00633614456663
$ws.Cells.Item(18, 7).Value = 3.620846481420698
$ws.Cells.Item(18, 10).Value = 9.569325923302925
$ws.Cells.Item(18, 13).Value = 22.06932933876651
$ws.Cells.Item(18, 14).Value = 18.11866569897416
$ws.Cells.Item(18, 15).Value = 24.97839217462604

# Row 19
$ws.Cells.Item(19, 2).Value = 8.02048300675448
$ws.Cells.Item(19, 4).Value = 8.923384823005231
$ws.Cells.Item(19, 5).Value = 12.85617917677946
$ws.Cells.Item(19, 6).Value = 32.988000281829
$ws.Cells.Item(19, 7).Value = 3.621082365936921
$ws.Cells.Item(19, 10).Value = 9.572469624465739
$ws.Cells.Item(19, 13).Value = 22.02392613724794
$ws.Cells.Item(19, 14).Value = 18.11087562111241
$ws.Cells.Item(19, 15).Value = 24.9642555615469

# Row 20
$ws.Cells.Item(20, 2).Value = 8.053786665256872
$ws.Cells.Item(20, 4).Value = 8.914508971826377
$ws.Cells.Item(20, 5).Value = 12.83427739262034
$ws.Cells.Item(20, 6).Value = 33.07109877202829
$ws.Cells.Item(20, 7).Value = 3.620027328819245
$ws.Cells.Item(20, 10).Value = 9.558404533397177
$ws.Cells.Item(20, 13).Value = 22.22756541615401
$ws.Cells.Item(20, 14).Value = 18.14606084930414
$ws.Cells.Item(20, 15).Value = 25.02829081668965

# Row 21
$ws.Cells.Item(21, 2).Value = 8.166580817353216
$ws.Cells.Item(21, 4).Value = 8.885817331839615
$ws.Cells.Item(21, 5).Value = 12.76308156047631
$ws.Cells.Item(21, 6).Value = 33.36034343387343
$ws.Cells.Item(21, 7).Value = 3.616595418789633
$ws.Cells.Item(21, 10).Value = 9.51257829684922
$ws.Cells.Item(21, 13).Value = 22.89931227982074
$ws.Cells.Item(21, 14).Value = 18.26646205478932
$ws.Cells.Item(21, 15).Value = 25.2506230968352

# Row 22
$ws.Cells.Item(22, 2).Value = 8.240847049667297
$ws.Cells.Item(22, 4).Value = 8.867906219286539
$ws.Cells.Item(22, 5).Value = 12.71832938113725
$ws.Cells.Item(22, 6).Value = 33.55662562667731
$ws.Cells.Item(22, 7).Value = 3.614436173422911
$ws.Cells.Item(22, 10).Value = 9.483690850014341
$ws.Cells.Item(22, 13).Value = 23.32842420788545
$ws.Cells.Item(22, 14).Value = 18.34668807776555
$ws.Cells.Item(22, 15).Value = 25.4010976784792

# Row 23
$ws.Cells.Item(23, 2).Value = 8.201167284723942
$ws.Cells.Item(23, 4).Value = 8.877388914640278
$ws.Cells.Item(23, 5).Value = 12.74205197369525
$ws.Cells.Item(23, 6).Value = 33.4512192982494
$ws.Cells.Item(23, 7).Value = 3.615580970917289
$ws.Cells.Item(23, 10).Value = 9.499011618452878
$ws.Cells.Item(23, 13).Value = 23.10034491052511
$ws.Cells.Item(23, 14).Value = 18.30373518271591
$ws.Cells.Item(23, 15).Value = 25.32032594763097

# Row 24
$ws.Cells.Item(24, 2).Value = 8.051954098651519
$ws.Cells.Item(24, 4).Value = 8.91499213574275
$ws.Cells.Item(24, 5).Value = 12.83547112859143
$ws.Cells.Item(24, 6).Value = 33.06649700964625
$ws.Cells.Item(24, 7).Value = 3.620084840727941
$ws.Cells.Item(24, 10).Value = 9.559171530109516
$ws.Cells.Item(24, 13).Value = 22.21642768262278
$ws.Cells.Item(24, 14).Value = 18.1441202008676
$ws.Cells.Item(24, 15).Value = 25.0247468247318

# Row 25
$ws.Cells.Item(25, 2).Value = 7.894143449742988
$ws.Cells.Item(25, 4).Value = 8.95917306400753
$ws.Cells.Item(25, 5).Value = 12.9438990912984
$ws.Cells.Item(25, 6).Value = 32.68414801350109
$ws.Cells.Item(25, 7).Value = 3.625305100456666
$ws.Cells.Item(25, 10).Value = 9.628649418443038
$ws.Cells.Item(25, 13).Value = 21.22465932447412
$ws.Cells.Item(25, 14).Value = 17.97913044024589
$ws.Cells.Item(25, 15).Value = 24.72927997271959
